# ArfGAP1 has an amphipathic helix (AH) that is split into two sub-helices
# (AH1 / AH2). This adds those two rows to the AA_seq / Protein_Name / AH# /
# Category table, matching the layout of the existing ArfGAP1 entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new shared-string values in the same order the author's Excel
# session produced them (B52, B51, A51, A52) so the rebuilt sharedStrings.xml
# table lines up with the target workbook.
$ws.Range("B52").Value = "ArfGAP1-AH2"
$ws.Range("B51").Value = "ArfGAP1-AH1"
$ws.Range("A51").Value = "FLNNAMSSLYSGWSSFTTGASRFASAAKEGATKFGS"
$ws.Range("A52").Value = "IFDDVSSGVSQLASKVQGVGSKGWRDVTTFFS"

$ws.Range("C51").Value = 1
$ws.Range("D51").Value = 4
$ws.Range("C52").Value = 1
$ws.Range("D52").Value = 4

# Column A had to widen a lot to fit the newly-added sequences legibly.
$ws.Columns.Item(1).ColumnWidth = 126.33072916666667

# Selection / scroll moved down to show the newly added rows.
$ws.Range("A53").Select()
$excel.ActiveWindow.ScrollRow = 30
